$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty task numbers for "Мартин С. Цингилев" (row 5)
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 1

# Move the active selection as recorded after the edit
$ws.Range("C17").Select()
